$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1
$ws.Range("D8").Value = "server down 03.07.2019"
$ws.Range("D10").Value = "too slow"

$ws.Range("D10").Select()
